$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price for row 2 (meaningful value change)
$ws.Range("C2").Value = 12.15

# Refresh cached computed values (Wt., Complex., Freq., Power, Speed) for "group 1" rows
$group1 = @(2, 4, 6, 9, 11)
foreach ($r in $group1) {
    $ws.Range("D" + $r).Value = 1.899999999999995
    $ws.Range("E" + $r).Value = -9.8
    $ws.Range("F" + $r).Value = 20.7
    $ws.Range("G" + $r).Value = 16.1
    $ws.Range("H" + $r).Value = -8.9
}

# Refresh cached computed values for "group 2" rows
$group2 = @(3, 5, 7, 8, 10, 12)
foreach ($r in $group2) {
    $ws.Range("D" + $r).Value = -6.700000000000005
    $ws.Range("E" + $r).Value = 2.1
    $ws.Range("F" + $r).Value = 25.7
    $ws.Range("G" + $r).Value = 16.1
    $ws.Range("H" + $r).Value = 15.4
}
